# Insert one new data row at row 460 (pushing existing rows 460..553 down to
# 461..554) and populate it with the new weekly observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(460).Insert()

$ws.Range("A460").Value = 6
$ws.Range("B460").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C460").Value = "Metropolitana"
$ws.Range("D460").Value = 45015
$ws.Range("E460").Value = 13
$ws.Range("F460").Value = 100112032
$ws.Range("G460").Value = "Zapallo italiano"
$ws.Range("H460").Value = "Sin especificar"
$ws.Range("I460").Value = "Primera"
$ws.Range("J460").Value = 400
$ws.Range("K460").Value = 5000
$ws.Range("L460").Value = 6000
$ws.Range("M460").Value = 5425
$ws.Range("N460").Value = "$/caja 50 unidades"
$ws.Range("O460").Value = "Región Metropolitana"
$ws.Range("P460").Value = 108
$ws.Range("Q460").Value = 50
$ws.Range("R460").Value = "Hortaliza"
